$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.47"
$ws.Range("E2").Value = "'4.72%"
$ws.Range("D3").Value = "'34.82"
$ws.Range("E3").Value = "'12.29%"
$ws.Range("D4").Value = "'5.157"
$ws.Range("E4").Value = "'4.69%"
$ws.Range("D5").Value = "'0.07752"
$ws.Range("D6").Value = "'2.365"
$ws.Range("E6").Value = "'5.64%"
$ws.Range("D7").Value = "'8.016"
$ws.Range("E7").Value = "'3.83%"
$ws.Range("D8").Value = "'3.948"
$ws.Range("E8").Value = "'5.24%"
$ws.Range("D9").Value = "'0.9292"
$ws.Range("E9").Value = "'1.88%"
$ws.Range("D10").Value = "'0.09918"
$ws.Range("E10").Value = "'11.74%"
$ws.Range("D11").Value = "'0.1803"
$ws.Range("E11").Value = "'6.99%"
$ws.Range("D12").Value = "'0.08646"
$ws.Range("E12").Value = "'5.18%"
$ws.Range("D13").Value = "'0.03317"
$ws.Range("E13").Value = "'6.24%"
$ws.Range("D14").Value = "'0.09888"
$ws.Range("E14").Value = "'-0.94%"
$ws.Range("D15").Value = "'0.001495"
$ws.Range("E15").Value = "'-0.45%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005756"
$ws.Range("E16").Value = "'-1.31%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.463"
$ws.Range("E17").Value = "'-0.92%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.136"
$ws.Range("E18").Value = "'2.27%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3368"
$ws.Range("E19").Value = "'1.22%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "'0.1334"
$ws.Range("E20").Value = "'2.79%"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").Value = "'4.355"
$ws.Range("E21").Value = "'9.52%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.2301"
$ws.Range("E22").Value = "'5.11%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.04577"
$ws.Range("E23").Value = "'0.51%"
$ws.Range("D24").Value = "'0.001217"
$ws.Range("E24").Value = "'0.42%"
$ws.Range("D25").Value = "'0.004462"
$ws.Range("E25").Value = "'-2.54%"
$ws.Range("D26").Value = "'0.0001300"
$ws.Range("E26").Value = "'-0.16%"
$ws.Range("D39").Value = "'0.01789"
$ws.Range("E39").Value = "'12.73%"
$ws.Range("D40").Value = "'0.04791"
$ws.Range("E40").Value = "'7.38%"
$ws.Range("D41").Value = "'0.007739"
$ws.Range("E41").Value = "'5.98%"
$ws.Range("D42").Value = "'0.1411"
$ws.Range("E42").Value = "'6.64%"
$ws.Range("D43").Value = "'0.007165"
$ws.Range("E43").Value = "'-25.04%"
$ws.Range("D44").Value = "'0.002099"
$ws.Range("E44").Value = "'-7.28%"
$ws.Range("D45").Value = "'0.009179"
$ws.Range("E45").Value = "'3.00%"
$ws.Range("D46").Value = "'0.00006123"
$ws.Range("E46").Value = "'0.27%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.15%"
$ws.Range("E48").Value = "'45.07%"
$ws.Range("D49").Value = "'0.002001"
$ws.Range("E49").Value = "'-0.15%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.15%"
$ws.Range("E51").Value = "'-0.15%"
